$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 5 with the missing X5/Y5 values
$ws.Range("X5").Value = -0.59999799999999937
$ws.Range("Y5").Value = "Down"

# Add new row 6 (scan results)
$ws.Range("A6").Value = 42647.883194444446
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 14
$ws.Range("E6").Value = 19796
$ws.Range("F6").Value = 2977
$ws.Range("G6").Value = 54
$ws.Range("H6").Value = 41
$ws.Range("I6").Value = 78
$ws.Range("J6").Value = 21
$ws.Range("K6").Value = 17251
$ws.Range("L6").Value = 387
$ws.Range("M6").Value = 297
$ws.Range("N6").Value = 84
$ws.Range("O6").Value = 23
$ws.Range("P6").Value = "Noun"
$ws.Range("Q6").Value = 53.235658945584888
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = -0.0862
$ws.Range("T6").Value = -0.0166
$ws.Range("U6").Value = 6.69
$ws.Range("V6").Value = 1.88
$ws.Range("W6").Value = 0

# Match the percentage formatting used by the rest of the table (column A
# already inherits the date style from the column definition)
$ws.Range("S6:T6").NumberFormat = "0.00%"
